# Apply "16th May Refresh" update to master-reg_center_user sheet:
# append 3 new rows (10005 / 110033-110035) with same lang/active/cr_by/cr_dtimes
# values as the existing rows, then leave the selection on the first empty row
# (matching Excel's behaviour after entering data at the bottom of a table).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @(10005, 110033, "eng", $true, "superadmin", "now()"),
    @(10005, 110034, "eng", $true, "superadmin", "now()"),
    @(10005, 110035, "eng", $true, "superadmin", "now()")
)

$startRow = 34
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Scroll back to top-left and move the selection to the row right after the
# newly entered data, selecting the rest of the sheet below it (as Excel does
# after typing a block of rows and pressing Ctrl+Down/selecting to the end).
$ws.Application.Goto($ws.Range("A1"), $true)
$ws.Range("A37:XFD1048576").Select()
